$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A (shifts all data right by one column)
$ws.Columns("A:A").Insert()

# Populate the new column A with the table name reference
$ws.Range("A1").Value = "Table"
$ws.Range("A2:A11").Value = "[Vend].[tblProject]"

# Update the generated INSERT statement formulas to reference the table name
# cell (A) instead of a hardcoded table name, and shift the column refs
$ws.Range("G2").Formula = '="INSERT INTO "&A2&" ([" & B$1 &"],["&C$1&"],["&D$1&"]) VALUES ( ''" & B2 & "'',''" & C2 & "'',''" & D2 & "'' )"'
$ws.Range("G3:G11").Formula = '="INSERT INTO "&A3&" ([" & B$1 &"],["&C$1&"],["&D$1&"]) VALUES ( ''" & B3 & "'',''" & C3 & "'',''" & D3 & "'' )"'

# Restore the active selection
$ws.Range("G6").Select()
